# Refresh the crypto price / 1h-volume table (GitHub Actions data pull).
# Price cells (column D) are numeric-looking text (e.g. "292.95") that must
# stay plain text, matching the original inline-string cells. Excel's COM
# Value setter auto-detects such strings as numbers, so each is written with
# a leading apostrophe (forces text) and then the cell style is reset to
# "Normal" so no stray quote-prefix/number-format style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '45.145.53'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.54%  '
$ws.Range('D3').Value = "'" + '2.386.77'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.17%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'" + '292.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.17%  '
$ws.Range('D6').Value = "'" + '93.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.75%  '
$ws.Range('D7').Value = "'" + '0.556'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.18%  '
$ws.Range('D8').Value = "'" + '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'" + '0.499'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.27%  '
$ws.Range('D10').Value = "'" + '34.14'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.92%  '
$ws.Range('D11').Value = "'" + '0.0776'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('D12').Value = "'" + '6.97'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.01%  '
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').Value = "'" + '2.750.69'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.11%  '
$ws.Range('D15').Value = "'" + '2.391.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.34%  '
$ws.Range('D16').Value = "'" + '13.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.18%  '
$ws.Range('D17').Value = "'" + '0.822'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.93%  '
$ws.Range('D18').Value = "'" + '45.108.40'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.60%  '
$ws.Range('D19').Value = "'" + '12.42'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.08%  '
$ws.Range('D20').Value = "'" + '0.0₃0929'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.29%  '
$ws.Range('D21').Value = "'" + '6.08'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.66%  '
$ws.Range('D22').Value = "'" + '66.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.46%  '
$ws.Range('D23').Value = "'" + '237.84'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.07%  '
$ws.Range('D24').Value = "'" + '2.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.36%  '
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = "'" + '1.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('D28').Value = "'" + '37.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -13.40%  '
$ws.Range('D29').Value = "'" + '9.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.38%  '
$ws.Range('E30').Value = '  +16.12%  '
$ws.Range('D31').Value = "'" + '20.90'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.07%  '
$ws.Range('D32').Value = "'" + '147.35'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('D33').Value = "'" + '2.69'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.57%  '
$ws.Range('D34').Value = "'" + '5.40'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('D35').Value = "'" + '0.0757'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.39%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = "'" + '1.95'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +12.08%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = "'" + '0.111'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.03%  '
$ws.Range('D38').Value = "'" + '0.114'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.50%  '
$ws.Range('E39').Value = '  -12.23%  '
$ws.Range('D40').Value = "'" + '3.69'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.92%  '
$ws.Range('D41').Value = "'" + '0.0292'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.37%  '
$ws.Range('D42').Value = "'" + '1.967.83'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.78%  '
$ws.Range('D43').Value = "'" + '3.15'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.52%  '
$ws.Range('D44').Value = "'" + '0.998'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = "'" + '88.15'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.62%  '
$ws.Range('E46').Value = '  -14.64%  '
$ws.Range('D47').Value = "'" + '8.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.63%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = "'" + '99.25'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.42%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = "'" + '2.619.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.07%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'" + '14.72'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +14.69%  '
$ws.Range('D51').Value = "'" + '0.182'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.38%  '
